$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a serial date value that was bumped by one
# day (45179 -> 45180) for every data row (rows 2 through 292).
$ws.Range("C2:C292").Value = 45180
